$d = $word.ActiveDocument

# Locate the paragraph that ends the existing "Bar03" bullet item -- the new
# "MSMap05/06/07" bullet items need to be inserted right after it (and before
# the trailing blank paragraph that precedes the sectPr).
$anchor = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs.Item($i)
    if ($candidate.Range.Text -like "*desired Bar03 barplot*") {
        $anchor = $candidate
    }
}
if ($anchor -eq $null) {
    throw "Could not find the 'desired Bar03 barplot' paragraph to anchor the insertion."
}

# Collapsed insertion point right after the Bar03 paragraph's own end-of-
# paragraph mark (i.e. right before the paragraph that currently follows it).
$insPoint = $d.Range($anchor.Range.End, $anchor.Range.End)

function New-BulletParagraphXml($boldText, $restText) {
    return "<w:p><w:pPr><w:pStyle w:val=`"Normal`"/><w:numPr><w:ilvl w:val=`"0`"/><w:numId w:val=`"1`"/></w:numPr><w:rPr><w:rFonts w:eastAsia=`"Calibri`" w:cs=`"`" w:cstheme=`"minorBidi`" w:eastAsiaTheme=`"minorHAnsi`"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:eastAsia=`"Calibri`" w:cs=`"`" w:cstheme=`"minorBidi`" w:eastAsiaTheme=`"minorHAnsi`"/><w:b/><w:bCs/><w:color w:val=`"000000`"/><w:kern w:val=`"0`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/><w:lang w:val=`"en-US`" w:eastAsia=`"en-US`" w:bidi=`"ar-SA`"/></w:rPr><w:t>$boldText</w:t></w:r><w:r><w:rPr><w:rFonts w:eastAsia=`"Calibri`" w:cs=`"`" w:cstheme=`"minorBidi`" w:eastAsiaTheme=`"minorHAnsi`"/><w:b w:val=`"false`"/><w:bCs w:val=`"false`"/><w:color w:val=`"000000`"/><w:kern w:val=`"0`"/><w:sz w:val=`"22`"/><w:szCs w:val=`"22`"/><w:lang w:val=`"en-US`" w:eastAsia=`"en-US`" w:bidi=`"ar-SA`"/></w:rPr><w:t>$restText</w:t></w:r></w:p>"
}

$bodyXml = ""
$bodyXml += New-BulletParagraphXml "MSMap05_Data1.txt." ".. for one of three supplemental material maps"
$bodyXml += New-BulletParagraphXml "MSMap06_Data1.txt." ".. for one of three supplemental material maps"
$bodyXml += New-BulletParagraphXml "MSMap07_Data1.txt." ".. for one of three supplemental material maps"
# A trailing empty paragraph absorbs the merge that InsertXML always performs
# between the last inserted paragraph and whatever paragraph originally sat
# at the insertion point -- it gets deleted again right after the insert.
$bodyXml += "<w:p/>"

$packageXml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' + $bodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$anchorIndex = $anchor.Index
$insPoint.InsertXML($packageXml)

# The synthetic trailing empty paragraph (inserted solely to take the merge
# with the original next paragraph) is now the 4th new paragraph after the
# anchor; remove it so the original trailing paragraph is restored intact.
$junk = $d.Paragraphs.Item($anchorIndex + 4)
$junk.Range.Delete()
